# V1000_GRP_DESC.xlsx update
#
# The parameter-group reference table drops four groups that are no longer
# used by the programmer (Motor 2 V/f Pattern, Motor 2 Setup, PM Motor Setup,
# and Auto-Tuning) and refreshes the LIST_IDX running offsets for every group
# that follows "Motor Setup" to reflect the new parameter counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-up so the row numbers used for the earlier deletions stay
# valid (deleting top-down would shift row 54 up before we got to it).

# Row 54: T1 / "Auto-Tuning"
$ws.Range("A54:D54").EntireRow.Delete() | Out-Null

# Rows 24-26: E3 / "Motor 2 V/f Pattern", E4 / "Motor 2 Setup",
# E5 / "PM Motor Setup"
$ws.Range("A24:D26").EntireRow.Delete() | Out-Null

# After the deletions, the rows from "Fault Detection During PG Speed
# Control" (F1) onward need their LIST_IDX (column D) values refreshed to
# the recalculated cumulative parameter offsets.
$listIdxUpdates = @{
    24 = 209
    25 = 217
    26 = 262
    27 = 298
    28 = 305
    29 = 310
    30 = 321
    31 = 324
    32 = 335
    33 = 342
    34 = 350
    35 = 359
    36 = 373
    37 = 381
    38 = 385
    39 = 395
    40 = 401
    41 = 420
    42 = 424
    43 = 427
    44 = 434
    45 = 435
    46 = 451
    47 = 456
    48 = 465
    49 = 474
    50 = 524
}
foreach ($row in $listIdxUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $listIdxUpdates[$row]
}

# Leave the selection where the edit ended up.
$ws.Range("D51").Select() | Out-Null
